$d = $word.ActiveDocument
$d.Content.Find.Execute("63×75=4725", $true, $false, $false, $false, $false, $true, 1, $false, "56×65=3640", 2) | Out-Null
$d.Content.Find.Execute("58×53=3074", $true, $false, $false, $false, $false, $true, 1, $false, "84×80=6720", 2) | Out-Null
$d.Content.Find.Execute("77×15=1155", $true, $false, $false, $false, $false, $true, 1, $false, "58×40=2320", 2) | Out-Null
$d.Content.Find.Execute("47×95=4465", $true, $false, $false, $false, $false, $true, 1, $false, "19×65=1235", 2) | Out-Null
$d.Content.Find.Execute("89×93=8277", $true, $false, $false, $false, $false, $true, 1, $false, "45×62=2790", 2) | Out-Null
$d.Content.Find.Execute("48×22=1056", $true, $false, $false, $false, $false, $true, 1, $false, "49×14=686", 2) | Out-Null
$d.Content.Find.Execute("11×46=506", $true, $false, $false, $false, $false, $true, 1, $false, "15×70=1050", 2) | Out-Null
$d.Content.Find.Execute("20×70=1400", $true, $false, $false, $false, $false, $true, 1, $false, "95×11=1045", 2) | Out-Null
$d.Content.Find.Execute("91×21=1911", $true, $false, $false, $false, $false, $true, 1, $false, "70×99=6930", 2) | Out-Null
$d.Content.Find.Execute("25×64=1600", $true, $false, $false, $false, $false, $true, 1, $false, "70×92=6440", 2) | Out-Null
$d.Content.Find.Execute("96×17=1632", $true, $false, $false, $false, $false, $true, 1, $false, "98×71=6958", 2) | Out-Null
$d.Content.Find.Execute("35×63=2205", $true, $false, $false, $false, $false, $true, 1, $false, "91×98=8918", 2) | Out-Null
$d.Content.Find.Execute("99×18=1782", $true, $false, $false, $false, $false, $true, 1, $false, "46×41=1886", 2) | Out-Null
$d.Content.Find.Execute("56×24=1344", $true, $false, $false, $false, $false, $true, 1, $false, "38×13=494", 2) | Out-Null
$d.Content.Find.Execute("41×84=3444", $true, $false, $false, $false, $false, $true, 1, $false, "66×93=6138", 2) | Out-Null
$d.Content.Find.Execute("23×50=1150", $true, $false, $false, $false, $false, $true, 1, $false, "73×36=2628", 2) | Out-Null
$d.Content.Find.Execute("94×84=7896", $true, $false, $false, $false, $false, $true, 1, $false, "85×41=3485", 2) | Out-Null
$d.Content.Find.Execute("62×95=5890", $true, $false, $false, $false, $false, $true, 1, $false, "55×26=1430", 2) | Out-Null
$d.Content.Find.Execute("12×87=1044", $true, $false, $false, $false, $false, $true, 1, $false, "90×67=6030", 2) | Out-Null
$d.Content.Find.Execute("99×11=1089", $true, $false, $false, $false, $false, $true, 1, $false, "30×60=1800", 2) | Out-Null
$d.Content.Find.Execute("78×27=2106", $true, $false, $false, $false, $false, $true, 1, $false, "68×40=2720", 2) | Out-Null
$d.Content.Find.Execute("21×39=819", $true, $false, $false, $false, $false, $true, 1, $false, "51×87=4437", 2) | Out-Null
$d.Content.Find.Execute("28×77=2156", $true, $false, $false, $false, $false, $true, 1, $false, "26×92=2392", 2) | Out-Null
$d.Content.Find.Execute("17×42=714", $true, $false, $false, $false, $false, $true, 1, $false, "68×50=3400", 2) | Out-Null
$d.Content.Find.Execute("62×88=5456", $true, $false, $false, $false, $false, $true, 1, $false, "13×13=169", 2) | Out-Null
